# Automatische test-sync: 2025-06-24 19:40:50
#
# Adds the new "BTW-nummer toevoegen" mail-log entry (row 4) to the "Logs"
# sheet, extends its conditional formatting ranges, adds the matching
# aggregate row to the "Dashboard" sheet, and widens the bar chart's
# category/value series references to include the new row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Logs sheet: append the new row of data
# ---------------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A4").Value = "BTW-nummer toevoegen"
$logs.Range("B4").Value = "mailmind.test@zohomail.eu"
$logs.Range("C4").Value = "Mijn BTW-nummer is niet vermeld op de factuur."
$logs.Range("D4").Value = "Factuur / Administratie"
$logs.Range("E4").Value = "Beste klant,`nDank u wel voor uw bericht. Om uw BTW-nummer toe te voegen aan uw factuur, hebben we uw klantgegevens nodig. Kunt u alstublieft uw factuurnummer en het juiste BTW-nummer delen, zodat we dit zo snel mogelijk voor u in orde kunnen maken?`nMet vriendelijke groet,`n[Bedrijfsnaam] E-mailassistent"
$logs.Range("F4").Value = "2025-06-24 19:39:50"
$logs.Range("G4").Value = "Ja"

# The multi-line text in E4 would otherwise leave a stale explicit
# row height behind; AutoFit puts the row back to the sheet default.
$logs.Rows.Item(4).EntireRow.AutoFit()

# Extend the existing conditional-formatting blocks (Categorie / Beantwoord
# columns) so they keep covering the data range through row 4.
$logs.Range("D2:D3").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D4"))
$logs.Range("G2:G3").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G4"))

# ---------------------------------------------------------------------------
# 2. Dashboard sheet: append the matching category/count row
# ---------------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A4").Value = "Factuur / Administratie"
$dash.Range("B4").Value = 1

# ---------------------------------------------------------------------------
# 3. Chart: widen the category/value series references to row 4
# ---------------------------------------------------------------------------
$chart = $dash.ChartObjects().Item(1).Chart
$series = $chart.SeriesCollection().Item(1)
$series.Formula = "=SERIES('Dashboard'!B1,'Dashboard'!`$A`$2:`$A`$4,'Dashboard'!`$B`$2:`$B`$4,1)"
